$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.218.50"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.560.22"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.53%  "
$ws.Range("E9").Value = "  +4.37%  "
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "3.017.69"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "63.106.43"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.86%  "
$ws.Range("D17").Value = "2.561.07"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").Value = "2.683.51"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.37%  "
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  +8.02%  "
$ws.Range("D32").Value = "0.0₃0826"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "466.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("E36").Value = "  +2.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.85%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "151.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0555"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.615"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0985"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0242"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("E50").Value = "  -1.93%  "
